$wb = $excel.ActiveWorkbook

# --- Update the status text: "Ready for handoff" -> "In Translation" ---
# The status shows up in three places: the "Overview" sheet's per-locale
# columns (zh-cn / de-de, i.e. E2 and F2), and each locale sheet's own
# "Status" column (C2).
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("E2").Value = "In Translation"
$wsOverview.Range("F2").Value = "In Translation"

$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("C2").Value = "In Translation"

$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("C2").Value = "In Translation"

# --- Narrow the per-locale status columns to match the new (shorter) text ---
# "Overview" sheet: the "zh-cn" and "de-de" columns (E and F) hold the status.
$wsOverview.Columns.Item(5).ColumnWidth = 12.5
$wsOverview.Columns.Item(6).ColumnWidth = 12.5

# Each locale sheet ("zh-cn", "de-de"): the "Status" column (C) holds the status.
$wsZhCn.Columns.Item(3).ColumnWidth = 12.5
$wsDeDe.Columns.Item(3).ColumnWidth = 12.5
